$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.008.17"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.551.61"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.21"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3923"
$ws.Range("E7").Value = "  +3.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3185"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.34"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07138"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.064"
$ws.Range("E11").Value = "  -5.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.611"
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.43"
$ws.Range("E14").Value = "  -7.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.607"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "1.552.70"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001089"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06560"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.82"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.146"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.28"
$ws.Range("E22").Value = "  -4.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  -4.94%  "
$ws.Range("D24").Value = "22.013.03"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.363"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.357"
$ws.Range("E26").Value = "  -5.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.13"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.31"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.865"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "1.728.17"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.11"
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9668"
$ws.Range("E32").Value = "  -9.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.774"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08260"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.059"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.602"
$ws.Range("E36").Value = "  -13.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02229"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.049"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05963"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2030"
$ws.Range("E41").Value = "  -5.37%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.57"
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5736"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.741"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.98"
$ws.Range("E46").Value = "  -4.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5496"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.36"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.854"
$ws.Range("E49").Value = "  -5.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.124"
$ws.Range("E50").Value = "  -3.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06802"
$ws.Range("E51").Value = "  -2.99%  "
